$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 646 ("「この子カンガルーは抱っこが好き」" post) entirely.
# All subsequent rows shift up by one, matching the target diff which
# renumbers rows 647-849 down to 646-848 and updates dimension to A1:C848.
$ws.Rows.Item(646).Delete()
